$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ArrayForm")

# --- Rows 7-13: lookup table used by the new array-operator formulas ---
# Column A (write in row order so shared-string ids line up with the fixture)
$ws.Range("A7").Value = "a"
$ws.Range("A8").Value = "b"
$ws.Range("A9").Value = "c"
$ws.Range("A10").Value = "a"
$ws.Range("A11").Value = "b"
$ws.Range("A12").Value = "c"
$ws.Range("A13").Value = "a"

# Column B (write in an order that yields shared-string order z, y, x)
$ws.Range("B11").Value = "z"
$ws.Range("B8").Value = "y"
$ws.Range("B7").Value = "x"
$ws.Range("B9").Value = "x"
$ws.Range("B10").Value = "y"
$ws.Range("B12").Value = "x"
$ws.Range("B13").Value = "y"

# Column C
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 4
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 7

# Array formulas exercising the new comparison/arithmetic operators on arrays
$ws.Range("E7").FormulaArray = '=SUM((A7:A13="a")*(B7:B13="y")*C7:C13)'
$ws.Range("E8").FormulaArray = '=SUM((A7:A13<>"b")*(B7:B13<>"y")*C7:C13)'
$ws.Range("E9").FormulaArray = '=SUM((A7:A13>"b")*(B7:B13<"z")*(C7:C13+3.5))'

# --- Rows 16-17: array multiplication of two ranges ---
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 3

$ws.Range("A17").Value = 3
$ws.Range("B17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 4

$ws.Range("G16:H17").FormulaArray = "=A16:B17*D16:E17"

# Match the recorded selection in the fixture
$ws.Range("G21").Select() | Out-Null
